$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "B"=1.02; "C"=1.030481359198216; "D"=1.03517641479953; "E"=0.9926147277508489; "F"=1.040590344860708; "I"=1.036205589407848; "J"=1.035622365164361; "K"=1.037973769952293; "L"=0.9955398523336033; "M"=1.04337226164183; "N"=1.037093067710457 }
  3 = @{ "B"=1.02; "C"=1.031309086804573; "D"=1.03580093504479; "E"=0.9936372048519304; "F"=1.041678286629802; "I"=1.036400608681015; "J"=1.036092186306957; "K"=1.038408122685384; "L"=0.9963617723202692; "M"=1.044269928187567; "N"=1.037563556052948 }
  4 = @{ "B"=1.02; "C"=1.031845007539789; "D"=1.036205290973253; "E"=0.9942998659930995; "F"=1.04238301040096; "I"=1.036525793041142; "J"=1.036395859145871; "K"=1.038688736622173; "L"=0.9968940712668345; "M"=1.044850939572726; "N"=1.037867660142124 }
  5 = @{ "B"=1.02; "C"=1.032070385288781; "D"=1.036375340483163; "E"=0.9945786998346017; "F"=1.04267945575486; "I"=1.036578179126281; "J"=1.036523442567172; "K"=1.038806600069114; "L"=0.997117960005301; "M"=1.045095234594419; "N"=1.037995424746518 }
  6 = @{ "B"=1.02; "C"=1.032108231663358; "D"=1.036403895918776; "E"=0.9946255319796338; "F"=1.042729240738318; "I"=1.036586960807683; "J"=1.036544859628504; "K"=1.038826383577338; "L"=0.9971555583673453; "M"=1.045136255000984; "N"=1.038016872222534 }
  7 = @{ "B"=1.02; "C"=1.031848018748635; "D"=1.036207562955984; "E"=0.9943035907982488; "F"=1.042386970813937; "I"=1.036526493975868; "J"=1.036397564239845; "K"=1.038690311939639; "L"=0.9968970624462087; "M"=1.044854203707004; "N"=1.037869367657527 }
  8 = @{ "B"=1.02; "C"=1.030761025423068; "D"=1.035387421782374; "E"=0.9929600610674301; "F"=1.040957863591418; "I"=1.036271705085888; "J"=1.035781211765831; "K"=1.038120652271623; "L"=0.995817528259106; "M"=1.043675598453534; "N"=1.037252139892304 }
  9 = @{ "B"=1.02; "C"=1.028848154874785; "D"=1.033944204482079; "E"=0.9906006454969559; "F"=1.038445390426595; "I"=1.035815052707117; "J"=1.034692612748316; "K"=1.037113503145393; "L"=0.9939188001724441; "M"=1.041600011607559; "N"=1.03616199493932 }
  10 = @{ "B"=1.02; "C"=1.027574698292978; "D"=1.032983475401437; "E"=0.989033133672735; "F"=1.036774346200139; "I"=1.035505485202581; "J"=1.033965250597909; "K"=1.036439886397456; "L"=0.9926553831429383; "M"=1.040217180406309; "N"=1.03543359985123 }
  11 = @{ "B"=1.02; "C"=1.027023717878347; "D"=1.032567823172718; "E"=0.988355674866747; "F"=1.036051707599336; "I"=1.035370228363649; "J"=1.033649919962135; "K"=1.036147697187991; "L"=0.9921088820399291; "M"=1.039618619476584; "N"=1.035117821409784 }
  12 = @{ "B"=1.02; "C"=1.026819125756208; "D"=1.032413485444059; "E"=0.9881042295826724; "F"=1.035783428426947; "I"=1.035319806364773; "J"=1.033532736246816; "K"=1.036039089645341; "L"=0.9919059725120875; "M"=1.039396320069213; "N"=1.035000471280148 }
  13 = @{ "B"=1.02; "C"=1.026863008430249; "D"=1.032446588942987; "E"=0.9881581567098651; "F"=1.035840968831291; "I"=1.035330630263621; "J"=1.033557875065376; "K"=1.036062389731131; "L"=0.9919494934313052; "M"=1.039444002553885; "N"=1.035025645798714 }
  14 = @{ "B"=1.02; "C"=1.027006804865273; "D"=1.032555064457567; "E"=0.9883348863814464; "F"=1.03602952868; "I"=1.035366064166833; "J"=1.033640234654243; "K"=1.036138721188412; "L"=0.9920921077337197; "M"=1.039600243466144; "N"=1.035108122347645 }
  15 = @{ "B"=1.02; "C"=1.027095411448276; "D"=1.032621907003568; "E"=0.9884438009545853; "F"=1.036145725324348; "I"=1.035387872117449; "J"=1.033690971738324; "K"=1.036185741548187; "L"=0.9921799884222134; "M"=1.039696513009725; "N"=1.035158931484204 }
  16 = @{ "B"=1.02; "C"=1.027611274311754; "D"=1.033011068365755; "E"=0.9890781214508737; "F"=1.036822325083078; "I"=1.035514436260401; "J"=1.0339861701674; "K"=1.036459267397616; "L"=0.9926916645766087; "M"=1.040256909486462; "N"=1.03545454912891 }
  17 = @{ "B"=1.02; "C"=1.027934978702762; "D"=1.033255273692014; "E"=0.989476357848556; "F"=1.037246988849393; "I"=1.035593502460303; "J"=1.034171239782192; "K"=1.036630707423074; "L"=0.9930127773699352; "M"=1.040608489139238; "N"=1.035639881563782 }
  18 = @{ "B"=1.02; "C"=1.028123831856036; "D"=1.033397748180236; "E"=0.9897087662937556; "F"=1.037494778174081; "I"=1.035639503455691; "J"=1.034279151221977; "K"=1.036730656341145; "L"=0.9932001317071769; "M"=1.040813580219633; "N"=1.035747946250189 }
  19 = @{ "B"=1.02; "C"=1.028188232963482; "D"=1.033446333962464; "E"=0.9897880325774034; "F"=1.037579283201482; "I"=1.035655168748708; "J"=1.034315940009332; "K"=1.036764727973452; "L"=0.9932640239640975; "M"=1.040883514456696; "N"=1.035784787281842 }
  20 = @{ "B"=1.02; "C"=1.027900243946582; "D"=1.033229069280183; "E"=0.9894336180360679; "F"=1.037201417119026; "I"=1.035585031501548; "J"=1.034151387336399; "K"=1.036612318607468; "L"=0.9929783193494215; "M"=1.040570765832849; "N"=1.035620000925239 }
  21 = @{ "B"=1.02; "C"=1.026964458552582; "D"=1.03252311961937; "E"=0.9882828385668249; "F"=1.035973998607955; "I"=1.035355634774057; "J"=1.033615983343198; "K"=1.036116245552064; "L"=0.9920501090198102; "M"=1.039554233502103; "N"=1.035083836596955 }
  22 = @{ "B"=1.02; "C"=1.026376478702453; "D"=1.032079574290105; "E"=0.9875604150241495; "F"=1.035203087824089; "I"=1.035210353597751; "J"=1.033279030768876; "K"=1.035803908620484; "L"=0.9914670000341481; "M"=1.038915289197115; "N"=1.03474640551133 }
  23 = @{ "B"=1.02; "C"=1.026688140987544; "D"=1.032314675882199; "E"=0.9879432794643023; "F"=1.035611684682214; "I"=1.035287469287445; "J"=1.0334576859986; "K"=1.035969525352913; "L"=0.991776070289318; "M"=1.039253987456781; "N"=1.034925314451972 }
  24 = @{ "B"=1.02; "C"=1.027915938958542; "D"=1.033240909822122; "E"=0.9894529299347244; "F"=1.03722200874384; "I"=1.035588859524682; "J"=1.034160357915094; "K"=1.036620627873664; "L"=0.9929938892766442; "M"=1.040587811307513; "N"=1.035628984243185 }
  25 = @{ "B"=1.02; "C"=1.029342367724323; "D"=1.034317067570802; "E"=0.9912096547607049; "F"=1.039094233887538; "I"=1.035934015310905; "J"=1.034974333010446; "K"=1.037374264594101; "L"=0.9944092447426414; "M"=1.042136446192537; "N"=1.036444115276536 }
}

foreach ($r in $data.Keys) {
  foreach ($c in $data[$r].Keys) {
    $ws.Range("$c$r").Value = $data[$r][$c]
  }
}
Write-Host "Updated cells"